$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @("ECs","Vegfc","Vipr2","ECs",3,1,2.454345666666666,7.363036999999999,0.2857469401257222,0.3027613480760606,2,0.6666666666666666,0.077915,0.233745,0.01344156357222124,0.01898188492565723,0.1912303426183333,1.721073083565,0.003840885661267591,0.005746981069116634),
        @("ECs","Vegfc","Vipr2","FAPs",3,1,2.454345666666666,7.363036999999999,0.2857469401257222,0.3027613480760606,3,1,0.6131856666666667,1.839557,0.1057841766036689,0.1493861228611831,1.504969583845444,13.544726254609,0.0302275047782174,0.04522834394130779),
        @("ECs","Vegfc","Vipr2","Inflammatory-Mac",3,1,2.454345666666666,7.363036999999999,0.2857469401257222,0.3027613480760606,1,0.3333333333333333,0.029863,0.089589,0.005151837424850708,0.007275313219982055,0.07329412464366666,0.659647121793,0.00147212178017627,0.002202683638157352),
        @("ECs","Vegfc","Vipr2","MuSCs",3,1,2.454345666666666,7.363036999999999,0.2857469401257222,0.3027613480760606,2,1,5.075609,10.151218,0.8756224223992591,0.8243566789931777,12.45729895484433,74.743793729066,0.2502064279060609,0.2495833394274788),
        @("FAPs","Vegfc","Vipr2","ECs",3,1,4.686805000000001,14.060415,0.5456607868665887,0.5781514068052169,2,0.6666666666666666,0.077915,0.233745,0.01344156357222124,0.01898188492565723,0.365172411575,3.286551704175,0.007334534155535516,0.01097440347358346),
        @("FAPs","Vegfc","Vipr2","FAPs",3,1,4.686805000000001,14.060415,0.5456607868665887,0.5781514068052169,3,1,0.6131856666666667,1.839557,0.1057841766036689,0.1493861228611831,2.873881648461667,25.864934836155,0.05772227704359216,0.08636779708936997),
        @("FAPs","Vegfc","Vipr2","Inflammatory-Mac",3,1,4.686805000000001,14.060415,0.5456607868665887,0.5781514068052169,1,0.3333333333333333,0.029863,0.089589,0.005151837424850708,0.007275313219982055,0.139962057715,1.259658519435,0.002811155663052777,0.004206232573081218),
        @("FAPs","Vegfc","Vipr2","MuSCs",3,1,4.686805000000001,14.060415,0.5456607868665887,0.5781514068052169,3,1,5.075609,10.151218,0.8756224223992591,0.8243566789931777,23.788389639245,142.73033783547,0.4777928200044083,0.4766029736691823),
        @("MuSCs","Vegfc","Vipr2","ECs",2,1,1.4480775,2.896155,0.1685922730076891,0.1190872451187225,2,0.6666666666666666,0.029863,0.233745,0.01344156357222124,0.01898188492565723,0.1128269584125,0.6769617504749998,0.01854368572461079,0.002260500382957126),
        @("MuSCs","Vegfc","Vipr2","FAPs",2,1,1.4480775,2.896155,0.1685922730076891,0.1190872451187225,3,1,0.6131856666666667,1.839557,0.1057841766036689,0.1493861228611831,0.8879403672225,5.327642203335,0.01783439478185935,0.0177899818305053),
        @("MuSCs","Vegfc","Vipr2","Inflammatory-Mac",2,1,1.4480775,2.896155,0.1685922730076891,0.1190872451187225,1,0.3333333333333333,0.029863,0.089589,0.005151837424850708,0.007275313219982055,0.0432439383825,0.259463630295,0.0008685599816216608,0.0008663970087434852),
        @("MuSCs","Vegfc","Vipr2","MuSCs",2,1,1.4480775,2.896155,0.1685922730076891,0.1190872451187225,2,1,5.075609,10.151218,0.8756224223992591,0.8243566789931777,7.3498751916975,29.39950076679,0.14762317448879,0.09817036589651659)
)

$numRows = $data.Count
$numCols = 20

$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$r1 = $ws.Cells.Item($startRow, 1)
$r2 = $ws.Cells.Item($endRow, $numCols)
$range = $ws.Range($r1, $r2)
$range.Value2 = $arr
